$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the second (AND) table at columns K:Q, mirroring the structure of
# the original table at columns A:G but shifted right by 10 columns, and
# using a uniform ".3-.7" marker along the diagonal instead of the varying
# values used in the original table.

# Header: "Known Parents" merged across M1:Q1, centered (copy style from C1:G1)
$ws.Range("C1:G1").Copy()
$ws.Range("M1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("M1").Value = "Known Parents"
$ws.Range("M1:Q1").Merge()

# Row 2: numeric headers 0..4
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 4

# Row label column K and numeric column L (mirrors A and B)
$ws.Range("K3").Value = "Parents"
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 4

# Fill the 5x5 matrix M3:Q7 with "x" everywhere, then set the diagonal
# (and below-diagonal) entries to match the AND table pattern.
$letters = @("M", "N", "O", "P", "Q")
for ($r = 3; $r -le 7; $r++) {
    for ($i = 0; $i -lt 5; $i++) {
        $col = $letters[$i]
        if ($i -lt ($r - 3)) {
            $ws.Range("$col$r").Value = 0
        } elseif ($i -eq ($r - 3)) {
            $ws.Range("$col$r").Value = ".3-.7"
        } else {
            $ws.Range("$col$r").Value = "x"
        }
    }
}

# Update selection / active cell to match the authored state
$ws.Range("N10").Select()
